$d = $word.ActiveDocument

# wdReplace constants
$wdReplaceNone = 0
$wdFindWrap   = 1

# 1) Locataire header line: drop the stray spaces inside the "${ first_name}" /
#    "${ last_name}" placeholders.
#    "${gender} ${ first_name} ${ last_name}" -> "${gender} ${first_name} ${last_name}"
$d.Content.Find.Execute(
    '${gender} ${ first_name} ${ last_name}', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${gender} ${first_name} ${last_name}', 2) | Out-Null

# 2) "des locaux sis" / Locataire-address placeholders: merge the split
#    "${" + "street" + "}" runs back into a single "${street}" token
#    (text itself is unchanged, this just tidies the template token).
$d.Content.Find.Execute(
    '${street}', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${street}', 2) | Out-Null

# 3) Same cleanup for the "${cp} ${city}" placeholders.
$d.Content.Find.Execute(
    '${cp} ${city}', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${cp} ${city}', 2) | Out-Null

# 4) "Reçu pour : ${loca_first_name} ${loca_last_name}" -> use the same
#    ${first_name} / ${last_name} placeholders as the header above.
$d.Content.Find.Execute(
    '${loca_first_name} ${loca_last_name}', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${first_name} ${ last_name}', 2) | Out-Null

# 5) "${first_day} au ${last_day} ${date_2}" -> "...${month}"
$d.Content.Find.Execute(
    '${date_2}', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${month}', 2) | Out-Null

# 6) "Loyer nu : ${loyer_hc} €" placeholder cleanup (merge split runs).
$d.Content.Find.Execute(
    '${loyer_hc} ', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${loyer_hc} ', 2) | Out-Null

# 7) "Charges : ${charges} €" placeholder cleanup (merge split runs).
$d.Content.Find.Execute(
    '${charges} ', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${charges} ', 2) | Out-Null

# 8) Footer quittance number: "Le ${date_3}" -> "Le ${date_1}"
$d.Content.Find.Execute(
    '${date_3}', $true, $false, $false, $false, $false,
    $true, $wdFindWrap, $false, '${date_1}', 2) | Out-Null
